$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.407.72'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '1.878.40'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7165'
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.76'
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07977'
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3149'
$ws.Range("E9").Value = '  +0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.95'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08101'
$ws.Range("E11").Value = '  -3.43%  '
$ws.Range("D12").Value = '1.884.12'
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.70'
$ws.Range("E13").Value = '  +3.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.231'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7080'
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.408'
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008441'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").Value = '29.411.82'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.92'
$ws.Range("E19").Value = '  +5.04%  '
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").Value = '2.135.16'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.678'
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.075'
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.36'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("E28").Value = '  +2.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.421'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.317'
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.222'
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.948'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01886'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").Value = '1.277.08'
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.762'
$ws.Range("E40").Value = '  +1.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.411'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9080'
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '111.78'
$ws.Range("E43").Value = '  +0.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '74.24'
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000130'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").Value = '2.030.54'
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.807'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4345'
$ws.Range("E51").Value = '  -0.35%  '
